$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "68.192.40"
$ws.Range("E2").Value = "  +6.79%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.616.71"
$ws.Range("E3").Value = "  +3.45%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.07%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'417.04"
$ws.Range("E5").Value = "  +0.50%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'129.85"
$ws.Range("E6").Value = "  -0.16%  "

# Row 7 - XRP
$ws.Range("D7").Value = "'0.657"
$ws.Range("E7").Value = "  +3.83%  "

# Row 8 - LidoStakedEther
$ws.Range("D8").Value = "3.600.66"
$ws.Range("E8").Value = "  +3.19%  "

# Row 9 - USDC
$ws.Range("D9").Value = "'0.998"
$ws.Range("E9").Value = "  -0.12%  "

# Row 10 - Cardano
$ws.Range("D10").Value = "'0.758"
$ws.Range("E10").Value = "  +1.29%  "

# Row 11 - Dogecoin
$ws.Range("D11").Value = "'0.192"
$ws.Range("E11").Value = "  +22.08%  "

# Row 12 - ShibaInu
$ws.Range("D12").Value = "'0.0000418"
$ws.Range("E12").Value = "  +85.80%  "

# Row 13 - Avalanche
$ws.Range("D13").Value = "'41.89"
$ws.Range("E13").Value = "  -2.05%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "'9.78"
$ws.Range("E14").Value = "  +0.01%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "4.187.19"
$ws.Range("E15").Value = "  +3.75%  "

# Row 16 - TRON
$ws.Range("E16").Value = "  -0.20%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "3.618.45"
$ws.Range("E17").Value = "  +3.98%  "

# Row 18 - Chainlink
$ws.Range("E18").Value = "  -2.13%  "

# Row 19 - Polygon
$ws.Range("D19").Value = "'1.11"
$ws.Range("E19").Value = "  +1.13%  "

# Row 20 - WrappedBTC
$ws.Range("D20").Value = "67.906.36"
$ws.Range("E20").Value = "  +6.65%  "

# Row 21 - Uniswap
$ws.Range("D21").Value = "'12.27"
$ws.Range("E21").Value = "  -1.80%  "

# Row 22 - BitcoinCash
$ws.Range("D22").Value = "'460.20"
$ws.Range("E22").Value = "  +0.18%  "

# Row 23 - Litecoin
$ws.Range("D23").Value = "'88.86"
$ws.Range("E23").Value = "  -1.31%  "

# Rows 24/25 swap content: ImmutableX now at 24, InternetComputer(DFINITY) now at 25
$ws.Range("B24").Value = "ImmutableX"
$ws.Range("C24").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D24").Value = "'3.04"
$ws.Range("E24").Value = "  -7.00%  "

$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").Value = "'13.20"
$ws.Range("E25").Value = "  -0.90%  "

# Row 26 - Filecoin
$ws.Range("D26").Value = "'10.05"
$ws.Range("E26").Value = "  -2.22%  "

# Row 27 - PancakeSwap
$ws.Range("D27").Value = "'3.26"
$ws.Range("E27").Value = "  -1.86%  "

# Row 28 - EthereumClassic
$ws.Range("E28").Value = "  +4.82%  "

# Row 29 - LEO
$ws.Range("E29").Value = "  +3.87%  "

# Row 30 - Toncoin
$ws.Range("D30").Value = "'2.78"
$ws.Range("E30").Value = "  +4.40%  "

# Row 31 - Cosmos
$ws.Range("D31").Value = "'12.17"
$ws.Range("E31").Value = "  -4.38%  "

# Row 32 - Hedera
$ws.Range("E32").Value = "  +3.33%  "

# Row 33 - RenderToken
$ws.Range("E33").Value = "  -4.68%  "

# Row 34 - InjectiveProtocol
$ws.Range("D34").Value = "'40.17"
$ws.Range("E34").Value = "  -0.08%  "

# Row 35 - Kaspa
$ws.Range("E35").Value = "  -7.75%  "

# Row 36 - Dai
$ws.Range("E36").Value = "  -0.24%  "

# Row 37 - OKB
$ws.Range("D37").Value = "'55.86"
$ws.Range("E37").Value = "  -2.99%  "

# Row 38 - PEPE
$ws.Range("D38").Value = "0.0₃0789"
$ws.Range("E38").Value = "  +16.87%  "

# Row 39 - VeChain
$ws.Range("D39").Value = "'0.0483"
$ws.Range("E39").Value = "  -1.29%  "

# Row 40 - Stellar
$ws.Range("E40").Value = "  +7.91%  "

# Row 41 - FirstDigitalUSD
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  +0.19%  "

# Rows 42/43 swap content: Monero now at 42, WEMIXToken now at 43
$ws.Range("B42").Value = "Monero"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D42").Value = "'148.24"
$ws.Range("E42").Value = "  +2.01%  "

$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").Value = "'2.72"
$ws.Range("E43").Value = "  -2.88%  "

# Row 44 - Stacks
$ws.Range("D44").Value = "'2.94"
$ws.Range("E44").Value = "  -4.67%  "

# Row 45 - LidoDAOToken
$ws.Range("D45").Value = "'3.22"
$ws.Range("E45").Value = "  -3.47%  "

# Rows 46/47/48 rotate: NEARProtocol, Cronos, ThetaToken
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").Value = "'4.23"
$ws.Range("E46").Value = "  -8.52%  "

$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "'0.169"
$ws.Range("E47").Value = "  +19.97%  "

$ws.Range("B48").Value = "ThetaToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D48").Value = "'2.56"
$ws.Range("E48").Value = "  +9.49%  "

# Row 49 - TheGraph
$ws.Range("E49").Value = "  -4.70%  "

# Row 50 - ARBITRUM
$ws.Range("D50").Value = "'1.93"
$ws.Range("E50").Value = "  -4.00%  "

# Row 51 - ApeXProtocol
$ws.Range("E51").Value = "  +13.11%  "
